# Apply updated FFXIV Leve profit figures (scheduled data refresh)
$wb = $excel.ActiveWorkbook

# --- ALC sheet ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 375749.16
$ws.Range("J17").Value = 375749.16
$ws.Range("L17").Value = 1127247.48
$ws.Range("N17").Value = -1127583.48
$ws.Range("H33").Value = 58986.176
$ws.Range("I33").Value = 32433.484
$ws.Range("J33").Value = 333364
$ws.Range("K33").Value = 32433.484
$ws.Range("L33").Value = 333364
$ws.Range("M33").Value = -32204.484
$ws.Range("N33").Value = -333822
$ws.Range("H129").Value = 1003.6
$ws.Range("I129").Value = 450
$ws.Range("K129").Value = 1350
$ws.Range("M129").Value = 3650
$ws.Range("H132").Value = 3420.1
$ws.Range("I132").Value = 3013.152
$ws.Range("J132").Value = 8100
$ws.Range("K132").Value = 9039.456
$ws.Range("L132").Value = 24300
$ws.Range("M132").Value = -6509.456
$ws.Range("N132").Value = -29360

# --- ARM sheet ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6050.4424
$ws.Range("I32").Value = 5077.093
$ws.Range("K32").Value = 5077.093
$ws.Range("M32").Value = -4790.093
$ws.Range("H45").Value = 1748.6666
$ws.Range("I45").Value = 1815.5
$ws.Range("K45").Value = 1815.5
$ws.Range("M45").Value = -1438.5
$ws.Range("H74").Value = 1331.2444
$ws.Range("I74").Value = 809.8
$ws.Range("J74").Value = 5502.8
$ws.Range("K74").Value = 809.8
$ws.Range("L74").Value = 5502.8
$ws.Range("M74").Value = 64.20000000000005
$ws.Range("N74").Value = -7250.8
$ws.Range("H77").Value = 1331.2444
$ws.Range("I77").Value = 809.8
$ws.Range("J77").Value = 5502.8
$ws.Range("K77").Value = 4049
$ws.Range("L77").Value = 27514
$ws.Range("M77").Value = 319
$ws.Range("N77").Value = -36250
$ws.Range("H110").Value = 919.86664
$ws.Range("I110").Value = 674.1
$ws.Range("J110").Value = 1411.4
$ws.Range("K110").Value = 674.1
$ws.Range("L110").Value = 1411.4
$ws.Range("M110").Value = 1370.9
$ws.Range("N110").Value = -5501.4
$ws.Range("H132").Value = 5486.9355
$ws.Range("I132").Value = 7039.85
$ws.Range("K132").Value = 21119.55
$ws.Range("M132").Value = -18589.55

# --- BSM sheet ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1401.3513
$ws.Range("I20").Value = 1445.5883
$ws.Range("K20").Value = 1445.5883
$ws.Range("M20").Value = -1198.5883
$ws.Range("H130").Value = 35507.145
$ws.Range("J130").Value = 35507.145
$ws.Range("L130").Value = 35507.145
$ws.Range("N130").Value = -45547.145

# --- CRP sheet ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 750
$ws.Range("I16").Value = 750
$ws.Range("J16").Value = 750
$ws.Range("K16").Value = 750
$ws.Range("L16").Value = 750
$ws.Range("M16").Value = -463
$ws.Range("N16").Value = -1324
$ws.Range("H52").Value = 42750
$ws.Range("I52").Value = 0
$ws.Range("K52").Value = 0
$ws.Range("M52").ClearContents()
$ws.Range("H99").Value = 1540.6666
$ws.Range("I99").Value = 1564.2858
$ws.Range("J99").Value = 1520
$ws.Range("K99").Value = 1564.2858
$ws.Range("L99").Value = 1520
$ws.Range("M99").Value = -66.28580000000011
$ws.Range("N99").Value = -4516
$ws.Range("H107").Value = 1094.75
$ws.Range("I107").Value = 1557.0834
$ws.Range("J107").Value = 401.25
$ws.Range("K107").Value = 1557.0834
$ws.Range("L107").Value = 401.25
$ws.Range("M107").Value = 362.9166
$ws.Range("N107").Value = -4241.25
$ws.Range("H113").Value = 750
$ws.Range("I113").Value = 750
$ws.Range("J113").Value = 750
$ws.Range("K113").Value = 750
$ws.Range("L113").Value = 750
$ws.Range("M113").Value = 1420
$ws.Range("N113").Value = -5090
$ws.Range("H126").Value = 1540.6666
$ws.Range("I126").Value = 1564.2858
$ws.Range("J126").Value = 1520
$ws.Range("K126").Value = 4692.857400000001
$ws.Range("L126").Value = 4560
$ws.Range("M126").Value = -2222.857400000001
$ws.Range("N126").Value = -9500
$ws.Range("H132").Value = 7278
$ws.Range("I132").Value = 4875.385
$ws.Range("J132").Value = 15086.5
$ws.Range("K132").Value = 14626.155
$ws.Range("L132").Value = 45259.5
$ws.Range("M132").Value = -12096.155
$ws.Range("N132").Value = -50319.5

# --- CUL sheet ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 71
$ws.Range("I14").Value = 71
$ws.Range("K14").Value = 213
$ws.Range("M14").Value = -40
$ws.Range("H47").Value = 2285.4285
$ws.Range("I47").Value = 0
$ws.Range("J47").Value = 2285.4285
$ws.Range("K47").Value = 0
$ws.Range("M47").Value = 6856.2855
$ws.Range("N47").Value = -7718.2855
$ws.Range("H98").Value = 588594.75
$ws.Range("I98").Value = 379.57144
$ws.Range("K98").Value = 1138.71432
$ws.Range("M98").Value = 359.28568
$ws.Range("H113").Value = 2392
$ws.Range("I113").Value = 3447.2
$ws.Range("J113").Value = 633.3333
$ws.Range("K113").Value = 10341.6
$ws.Range("L113").Value = 1899.9999
$ws.Range("M113").Value = -8171.599999999999
$ws.Range("N113").Value = -6239.9999

# --- GSM sheet ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2358.8572
$ws.Range("I102").Value = 2128
$ws.Range("J102").Value = 2666.6667
$ws.Range("K102").Value = 2128
$ws.Range("L102").Value = 2666.6667
$ws.Range("M102").Value = -506
$ws.Range("N102").Value = -5910.6667
$ws.Range("H113").Value = 989.55554
$ws.Range("I113").Value = 984.7143
$ws.Range("K113").Value = 984.7143
$ws.Range("M113").Value = 1185.2857
$ws.Range("H122").Value = 1713.5
$ws.Range("I122").Value = 1600.1538
$ws.Range("J122").Value = 1877.2222
$ws.Range("K122").Value = 4800.4614
$ws.Range("L122").Value = 5631.6666
$ws.Range("M122").Value = -2350.4614
$ws.Range("N122").Value = -10531.6666
$ws.Range("H126").Value = 112191.336
$ws.Range("I126").Value = 200856.4
$ws.Range("J126").Value = 1360
$ws.Range("K126").Value = 602569.2
$ws.Range("L126").Value = 4080
$ws.Range("M126").Value = -600099.2
$ws.Range("N126").Value = -9020
$ws.Range("H132").Value = 5351.7417
$ws.Range("I132").Value = 5973.5654
$ws.Range("J132").Value = 3564
$ws.Range("K132").Value = 17920.6962
$ws.Range("L132").Value = 10692
$ws.Range("M132").Value = -15390.6962
$ws.Range("N132").Value = -15752

# --- LTW sheet ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1758.4117
$ws.Range("I7").Value = 1045.75
$ws.Range("J7").Value = 1977.6923
$ws.Range("K7").Value = 1045.75
$ws.Range("L7").Value = 1977.6923
$ws.Range("M7").Value = -933.75
$ws.Range("N7").Value = -2201.6923
$ws.Range("H16").Value = 1566.0834
$ws.Range("I16").Value = 1499.1
$ws.Range("J16").Value = 1901
$ws.Range("K16").Value = 1499.1
$ws.Range("L16").Value = 1901
$ws.Range("M16").Value = -1329.1
$ws.Range("N16").Value = -2241
$ws.Range("H40").Value = 2728.1155
$ws.Range("I40").Value = 2478.682
$ws.Range("J40").Value = 4100
$ws.Range("K40").Value = 2478.682
$ws.Range("L40").Value = 4100
$ws.Range("M40").Value = -2342.682
$ws.Range("N40").Value = -4372
$ws.Range("H61").Value = 1712.5
$ws.Range("I61").Value = 1557.1428
$ws.Range("K61").Value = 1557.1428
$ws.Range("M61").Value = -1355.1428
$ws.Range("H113").Value = 1712.5
$ws.Range("I113").Value = 1557.1428
$ws.Range("K113").Value = 1557.1428
$ws.Range("M113").Value = 612.8571999999999
$ws.Range("H126").Value = 1758.4117
$ws.Range("I126").Value = 1045.75
$ws.Range("J126").Value = 1977.6923
$ws.Range("K126").Value = 3137.25
$ws.Range("L126").Value = 5933.0769
$ws.Range("M126").Value = -667.25
$ws.Range("N126").Value = -10873.0769

# --- WVR sheet ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H74").Value = 7395.154
$ws.Range("J74").Value = 7395.154
$ws.Range("L74").Value = 7395.154
$ws.Range("N74").Value = -9267.154
$ws.Range("H77").Value = 7395.154
$ws.Range("J77").Value = 7395.154
$ws.Range("L77").Value = 22185.462
$ws.Range("N77").Value = -31545.462
$ws.Range("H126").Value = 1837777
$ws.Range("I126").Value = 1266889.4
$ws.Range("J126").Value = 4349683
$ws.Range("K126").Value = 3800668.2
$ws.Range("L126").Value = 13049049
$ws.Range("M126").Value = -3798198.2
$ws.Range("N126").Value = -13053989
$ws.Range("H132").Value = 3216.804
$ws.Range("I132").Value = 3513.6667
$ws.Range("J132").Value = 1831.4445
$ws.Range("K132").Value = 10541.0001
$ws.Range("L132").Value = 5494.333500000001
$ws.Range("M132").Value = -8011.000100000001
$ws.Range("N132").Value = -10554.3335

